$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T3").Value = 9.75
$ws.Range("AA3").Value = 7.2
$ws.Range("AH3").Value = 23

$ws.Range("G7").Value = 1.45
$ws.Range("H7").Value = 4.2
$ws.Range("I7").Value = 6.1
$ws.Range("L7").Value = 1.21
$ws.Range("M7").Value = 3.55
$ws.Range("O7").Value = 2.02
$ws.Range("R7").Value = 1.78
$ws.Range("S7").Value = 1.83
$ws.Range("T7").Value = 7.4
$ws.Range("U7").Value = 7.2
$ws.Range("V7").Value = 8
$ws.Range("W7").Value = 10
$ws.Range("X7").Value = 11.25
$ws.Range("Y7").Value = 24
$ws.Range("Z7").Value = 12.5
$ws.Range("AA7").Value = 8.25
$ws.Range("AB7").Value = 17.5
$ws.Range("AC7").Value = 75
$ws.Range("AD7").Value = 600
$ws.Range("AE7").Value = 17.5
$ws.Range("AF7").Value = 40
$ws.Range("AG7").Value = 19
$ws.Range("AH7").Value = 120
$ws.Range("AI7").Value = 65
$ws.Range("AJ7").Value = 60

$ws.Range("G13").Value = 3.1
$ws.Range("H13").Value = 3.1
$ws.Range("I13").Value = 2.25
$ws.Range("L13").Value = 1.45
$ws.Range("M13").Value = 2.37
$ws.Range("N13").Value = 2.32
$ws.Range("O13").Value = 1.47
$ws.Range("P13").Value = 1.5
$ws.Range("Q13").Value = 2.27
$ws.Range("R13").Value = 2.02
$ws.Range("S13").Value = 1.62
$ws.Range("T13").Value = 7.4
$ws.Range("U13").Value = 14.5
$ws.Range("V13").Value = 11.75
$ws.Range("W13").Value = 40
$ws.Range("X13").Value = 32
$ws.Range("Y13").Value = 50
$ws.Range("Z13").Value = 6.9
$ws.Range("AA13").Value = 6.1
$ws.Range("AB13").Value = 18.5
$ws.Range("AE13").Value = 6
$ws.Range("AF13").Value = 9.5
$ws.Range("AG13").Value = 9.75
$ws.Range("AH13").Value = 22
$ws.Range("AI13").Value = 22
$ws.Range("AJ13").Value = 45

$ws.Range("G14").Value = 3.45
$ws.Range("H14").Value = 3.1
$ws.Range("I14").Value = 2.1
$ws.Range("L14").Value = 1.5
$ws.Range("M14").Value = 2.25
$ws.Range("N14").Value = 2.42
$ws.Range("O14").Value = 1.44
$ws.Range("P14").Value = 1.53
$ws.Range("Q14").Value = 2.2
$ws.Range("R14").Value = 2.12
$ws.Range("S14").Value = 1.57
$ws.Range("T14").Value = 7.5
$ws.Range("U14").Value = 16
$ws.Range("V14").Value = 13.5
$ws.Range("Y14").Value = 65
$ws.Range("Z14").Value = 6.5
$ws.Range("AA14").Value = 6.2
$ws.Range("AB14").Value = 21
$ws.Range("AE14").Value = 5.5
$ws.Range("AG14").Value = 9.75
$ws.Range("AH14").Value = 19
$ws.Range("AI14").Value = 22
